$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (not auto-coerced to a number by
# Excel's input-parsing heuristics) and then drop back to the default
# "Normal" style so no stray per-cell style index is left behind.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# New FedEx tracking numbers (column P, "ShipmentTracking") for rows 2-26.
$newTracking = @{
    2  = "320018616270"
    3  = "320018616280"
    4  = "320018616317"
    5  = "320018616339"
    6  = "320018616372"
    7  = "320018616394"
    8  = "320018616420"
    9  = "320018616442"
    10 = "320018616475"
    11 = "320018616497"
    12 = "320018616534"
    13 = "320018616556"
    14 = "320018616589"
    15 = "320018616604"
    16 = "320018616637"
    17 = "320018616659"
    18 = "320018616692"
    19 = "320018616718"
    20 = "320018616740"
    21 = "320018616762"
    22 = "320018616795"
    23 = "320018616800"
    24 = "320018616810"
    25 = "320018616821"
    26 = "320018616832"
}

foreach ($row in $newTracking.Keys) {
    Set-TextValue $ws.Cells.Item($row, 16) $newTracking[$row]
}

# Rows 22-26 also got new ActualRate (column Q) values and flipped from
# PASS to FAIL (column R).
$newActualRate = @{
    22 = "`$276.24"
    23 = "`$489.85"
    24 = "`$354.26"
    25 = "`$132.19"
    26 = "`$1,382.91"
}

foreach ($row in $newActualRate.Keys) {
    Set-TextValue $ws.Cells.Item($row, 17) $newActualRate[$row]
    $ws.Cells.Item($row, 18).Value = "FAIL"
}
